# Auto-generated edit script: updates crypto price/volume table cells per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.830.68'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '2.240.41'
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''116.36'
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('D6').Value = '''265.74'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '''0.630'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').Value = '''0.608'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '''46.63'
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').Value = '''0.0930'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = '''9.15'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('E13').Value = '  -2.24%  '
$ws.Range('D14').Value = '''15.37'
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('E15').Value = '  +1.04%  '
$ws.Range('D16').Value = '2.578.20'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').Value = '2.264.12'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '43.095.87'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').Value = '''6.74'
$ws.Range('E20').Value = '  -1.67%  '
$ws.Range('D21').Value = '''71.52'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '''2.37'
$ws.Range('E22').Value = '  -6.16%  '
$ws.Range('D23').Value = '''231.52'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '''9.48'
$ws.Range('E25').Value = '  -2.24%  '
$ws.Range('D26').Value = '''12.15'
$ws.Range('E26').Value = '  +6.51%  '
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').Value = '''41.13'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').Value = '''2.25'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').Value = '''3.30'
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('D31').Value = '''172.79'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '''21.21'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').Value = '''0.0899'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').Value = '''5.60'
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('D35').Value = '''4.32'
$ws.Range('E35').Value = '  +9.54%  '
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').Value = '''0.0374'
$ws.Range('E37').Value = '  +4.19%  '
$ws.Range('D38').Value = '''4.66'
$ws.Range('E38').Value = '  -0.72%  '
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').Value = '''2.51'
$ws.Range('E40').Value = '  -6.55%  '
$ws.Range('D41').Value = '''13.44'
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('D42').Value = '''0.235'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').Value = '''71.11'
$ws.Range('E43').Value = '  -8.71%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('E45').Value = '  -3.64%  '
$ws.Range('D46').Value = '''5.65'
$ws.Range('E46').Value = '  -9.17%  '
$ws.Range('D47').Value = '''73.26'
$ws.Range('E47').Value = '  +30.53%  '
$ws.Range('D48').Value = '''0.652'
$ws.Range('E48').Value = '  +14.86%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''8.43'
$ws.Range('E49').Value = '  -3.45%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0992'
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '''1.24'
$ws.Range('E51').Value = '  -0.98%  '
